$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Psap"
$ws.Cells.Item(2,3).Value = "Gpr37l1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 107.565699
$ws.Cells.Item(2,8).Value = 322.697097
$ws.Cells.Item(2,9).Value = 0.02858134952685079
$ws.Cells.Item(2,10).Value = 0.02881178879290993
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.1139563333333333
$ws.Cells.Item(2,14).Value = 0.341869
$ws.Cells.Item(2,15).Value = 0.05897041328053401
$ws.Cells.Item(2,16).Value = 0.08592218649695475
$ws.Cells.Item(2,17).Value = 12.257792650477
$ws.Cells.Item(2,18).Value = 110.320133854293
$ws.Cells.Item(2,19).Value = 0.001685453993713786
$ws.Cells.Item(2,20).Value = 0.002475571889975278

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Psap"
$ws.Cells.Item(3,3).Value = "Gpr37l1"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 107.565699
$ws.Cells.Item(3,8).Value = 322.697097
$ws.Cells.Item(3,9).Value = 0.02858134952685079
$ws.Cells.Item(3,10).Value = 0.02881178879290993
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.818476
$ws.Cells.Item(3,14).Value = 3.636952
$ws.Cells.Item(3,15).Value = 0.941029586719466
$ws.Cells.Item(3,16).Value = 0.9140778135030453
$ws.Cells.Item(3,17).Value = 195.605642054724
$ws.Cells.Item(3,18).Value = 1173.633852328344
$ws.Cells.Item(3,19).Value = 0.02689589553313701
$ws.Cells.Item(3,20).Value = 0.02633621690293465

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Psap"
$ws.Cells.Item(4,3).Value = "Gpr37l1"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 176.0123163333334
$ws.Cells.Item(4,8).Value = 528.036949
$ws.Cells.Item(4,9).Value = 0.04676834326297298
$ws.Cells.Item(4,10).Value = 0.04714541652489843
$ws.Cells.Item(4,11).Value = 1
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.1139563333333333
$ws.Cells.Item(4,14).Value = 0.341869
$ws.Cells.Item(4,15).Value = 0.05897041328053401
$ws.Cells.Item(4,16).Value = 0.08592218649695475
$ws.Cells.Item(4,17).Value = 20.05771819085345
$ws.Cells.Item(4,18).Value = 180.519463717681
$ws.Cells.Item(4,19).Value = 0.002757948530663395
$ws.Cells.Item(4,20).Value = 0.004050837271128936

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Psap"
$ws.Cells.Item(5,3).Value = "Gpr37l1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 176.0123163333334
$ws.Cells.Item(5,8).Value = 528.036949
$ws.Cells.Item(5,9).Value = 0.04676834326297298
$ws.Cells.Item(5,10).Value = 0.04714541652489843
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.818476
$ws.Cells.Item(5,14).Value = 3.636952
$ws.Cells.Item(5,15).Value = 0.941029586719466
$ws.Cells.Item(5,16).Value = 0.9140778135030453
$ws.Cells.Item(5,17).Value = 320.0741729565747
$ws.Cells.Item(5,18).Value = 1920.445037739448
$ws.Cells.Item(5,19).Value = 0.04401039473230959
$ws.Cells.Item(5,20).Value = 0.0430945792537695

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Psap"
$ws.Cells.Item(6,3).Value = "Gpr37l1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1855.177368
$ws.Cells.Item(6,8).Value = 5565.532104
$ws.Cells.Item(6,9).Value = 0.4929403451290834
$ws.Cells.Item(6,10).Value = 0.4969147134924725
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.1139563333333333
$ws.Cells.Item(6,14).Value = 0.341869
$ws.Cells.Item(6,15).Value = 0.05897041328053401
$ws.Cells.Item(6,16).Value = 0.08592218649695475
$ws.Cells.Item(6,17).Value = 211.409210540264
$ws.Cells.Item(6,18).Value = 1902.682894862376
$ws.Cells.Item(6,19).Value = 0.02906889587491112
$ws.Cells.Item(6,20).Value = 0.04269599868578106

# Row 7
$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "Psap"
$ws.Cells.Item(7,3).Value = "Gpr37l1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1855.177368
$ws.Cells.Item(7,8).Value = 5565.532104
$ws.Cells.Item(7,9).Value = 0.4929403451290834
$ws.Cells.Item(7,10).Value = 0.4969147134924725
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.818476
$ws.Cells.Item(7,14).Value = 3.636952
$ws.Cells.Item(7,15).Value = 0.941029586719466
$ws.Cells.Item(7,16).Value = 0.9140778135030453
$ws.Cells.Item(7,17).Value = 3373.595519451168
$ws.Cells.Item(7,18).Value = 20241.57311670701
$ws.Cells.Item(7,19).Value = 0.4638714492541722
$ws.Cells.Item(7,20).Value = 0.4542187148066915

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Psap"
$ws.Cells.Item(8,3).Value = "Gpr37l1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1534.435017666667
$ws.Cells.Item(8,8).Value = 4603.305053
$ws.Cells.Item(8,9).Value = 0.4077156935146256
$ws.Cells.Item(8,10).Value = 0.4110029317566841
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.1139563333333333
$ws.Cells.Item(8,14).Value = 0.341869
$ws.Cells.Item(8,15).Value = 0.05897041328053401
$ws.Cells.Item(8,16).Value = 0.08592218649695475
$ws.Cells.Item(8,17).Value = 174.8585883515619
$ws.Cells.Item(8,18).Value = 1573.727295164057
$ws.Cells.Item(8,19).Value = 0.02404316294751701
$ws.Cells.Item(8,20).Value = 0.03531427055319297

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Psap"
$ws.Cells.Item(9,3).Value = "Gpr37l1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1534.435017666667
$ws.Cells.Item(9,8).Value = 4603.305053
$ws.Cells.Item(9,9).Value = 0.4077156935146256
$ws.Cells.Item(9,10).Value = 0.4110029317566841
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.818476
$ws.Cells.Item(9,14).Value = 3.636952
$ws.Cells.Item(9,15).Value = 0.941029586719466
$ws.Cells.Item(9,16).Value = 0.9140778135030453
$ws.Cells.Item(9,17).Value = 2790.333253186409
$ws.Cells.Item(9,18).Value = 16741.99951911846
$ws.Cells.Item(9,19).Value = 0.3836725305671086
$ws.Cells.Item(9,20).Value = 0.3756886612034911

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Psap"
$ws.Cells.Item(10,3).Value = "Gpr37l1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 90.30225350000001
$ws.Cells.Item(10,8).Value = 180.604507
$ws.Cells.Item(10,9).Value = 0.02399426856646732
$ws.Cells.Item(10,10).Value = 0.0161251494330351
$ws.Cells.Item(10,11).Value = 1
$ws.Cells.Item(10,12).Value = 0.3333333333333333
$ws.Cells.Item(10,13).Value = 0.1139563333333333
$ws.Cells.Item(10,14).Value = 0.341869
$ws.Cells.Item(10,15).Value = 0.05897041328053401
$ws.Cells.Item(10,16).Value = 0.08592218649695475
$ws.Cells.Item(10,17).Value = 10.29051370059717
$ws.Cells.Item(10,18).Value = 61.743082203583
$ws.Cells.Item(10,19).Value = 0.001414951933728704
$ws.Cells.Item(10,20).Value = 0.001385508096876506

# Row 11
$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Psap"
$ws.Cells.Item(11,3).Value = "Gpr37l1"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 90.30225350000001
$ws.Cells.Item(11,8).Value = 180.604507
$ws.Cells.Item(11,9).Value = 0.02399426856646732
$ws.Cells.Item(11,10).Value = 0.0161251494330351
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 1.818476
$ws.Cells.Item(11,14).Value = 3.636952
$ws.Cells.Item(11,15).Value = 0.941029586719466
$ws.Cells.Item(11,16).Value = 0.9140778135030453
$ws.Cells.Item(11,17).Value = 164.212480735666
$ws.Cells.Item(11,18).Value = 656.8499229426641
$ws.Cells.Item(11,19).Value = 0.02257931663273862
$ws.Cells.Item(11,20).Value = 0.0147396413361586

